$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194 (shifts existing rows 194:202 down to 195:203,
# matching the weekly entry added to the "Feria Lagunitas de Puerto Montt - Pera" series).
$ws.Rows(194).Insert()

$ws.Range("A194").Value = 4
$ws.Range("B194").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value = "Los Lagos"
$ws.Range("D194").Value = 44568
$ws.Range("E194").Value = 10
$ws.Range("F194").Value = "Fruta"
$ws.Range("G194").Value = 100104
$ws.Range("H194").Value = "Frutos de pepita"
$ws.Range("I194").Value = 100104005
$ws.Range("J194").Value = "Pera"
$ws.Range("K194").Value = "Packham's Triumph"
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 400
$ws.Range("N194").Value = 14000
$ws.Range("O194").Value = 15000
$ws.Range("P194").Value = 14500
$ws.Range("Q194").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R194").Value = "Región de O'Higgins"
$ws.Range("S194").Value = 967
$ws.Range("T194").Value = 15
